$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear column I (I1 header "LV_2024" and data cells I2, I4, I5, I12 that held
# values copied from other year columns) so the LV_2024 column becomes blank,
# matching the removal of the LV_2024 data in this translation table.
# I1 loses its formatting entirely (fully cleared header cell), while the
# data cells only lose their content and keep their existing style.
$ws.Cells.Item(1, 9).Clear()
$ws.Cells.Item(2, 9).ClearContents()
$ws.Cells.Item(4, 9).ClearContents()
$ws.Cells.Item(5, 9).ClearContents()
$ws.Cells.Item(12, 9).ClearContents()

# Update the selected cell to I9 (matches the saved cursor position in the file)
$ws.Range("I9").Select()
